$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updating DAMSLTag (column I) and DialogAct (column J) for the affected rows.

$ws.Range("I3").Value = "ba"
$ws.Range("J3").Value = "Appreciation"
$ws.Range("I12").Value = "sv"
$ws.Range("J12").Value = "Statement-opinion"
$ws.Range("I14").Value = "b"
$ws.Range("J14").Value = "Acknowledge (Backchannel)"
$ws.Range("I47").Value = "sd"
$ws.Range("J47").Value = "Statement-non-opinion"
$ws.Range("I62").Value = "sd"
$ws.Range("J62").Value = "Statement-non-opinion"
$ws.Range("I69").Value = "ba"
$ws.Range("J69").Value = "Appreciation"
$ws.Range("I71").Value = "ba"
$ws.Range("J71").Value = "Appreciation"
$ws.Range("I77").Value = "b"
$ws.Range("J77").Value = "Acknowledge (Backchannel)"
$ws.Range("I79").Value = "aa"
$ws.Range("J79").Value = "Agree/Accept"
$ws.Range("I82").Value = "sv"
$ws.Range("J82").Value = "Statement-opinion"
$ws.Range("I102").Value = "ba"
$ws.Range("J102").Value = "Appreciation"
$ws.Range("I129").Value = "%"
$ws.Range("J129").Value = "Uninterpretable"
$ws.Range("I131").Value = "sd"
$ws.Range("J131").Value = "Statement-non-opinion"
$ws.Range("I132").Value = "sd"
$ws.Range("J132").Value = "Statement-non-opinion"
$ws.Range("I134").Value = "%"
$ws.Range("J134").Value = "Uninterpretable"
$ws.Range("I138").Value = "b"
$ws.Range("J138").Value = "Acknowledge (Backchannel)"
$ws.Range("I155").Value = "%"
$ws.Range("J155").Value = "Uninterpretable"
$ws.Range("I166").Value = "ba"
$ws.Range("J166").Value = "Appreciation"
$ws.Range("I170").Value = "sd"
$ws.Range("J170").Value = "Statement-non-opinion"
$ws.Range("I208").Value = "ba"
$ws.Range("J208").Value = "Appreciation"
$ws.Range("I210").Value = "ba"
$ws.Range("J210").Value = "Appreciation"
$ws.Range("I219").Value = "sv"
$ws.Range("J219").Value = "Statement-opinion"
$ws.Range("I230").Value = "sd"
$ws.Range("J230").Value = "Statement-non-opinion"
$ws.Range("I237").Value = "sd"
$ws.Range("J237").Value = "Statement-non-opinion"
$ws.Range("I244").Value = "sv"
$ws.Range("J244").Value = "Statement-opinion"
$ws.Range("I247").Value = "ba"
$ws.Range("J247").Value = "Appreciation"
$ws.Range("I256").Value = "sv"
$ws.Range("J256").Value = "Statement-opinion"
$ws.Range("I259").Value = "ba"
$ws.Range("J259").Value = "Appreciation"
$ws.Range("I260").Value = "aa"
$ws.Range("J260").Value = "Agree/Accept"
$ws.Range("I269").Value = "ba"
$ws.Range("J269").Value = "Appreciation"
$ws.Range("I273").Value = "%"
$ws.Range("J273").Value = "Uninterpretable"
$ws.Range("I275").Value = "b"
$ws.Range("J275").Value = "Acknowledge (Backchannel)"
$ws.Range("I280").Value = "b"
$ws.Range("J280").Value = "Acknowledge (Backchannel)"
$ws.Range("I318").Value = "sv"
$ws.Range("J318").Value = "Statement-opinion"
$ws.Range("I331").Value = "sd"
$ws.Range("J331").Value = "Statement-non-opinion"
$ws.Range("I340").Value = "b"
$ws.Range("J340").Value = "Acknowledge (Backchannel)"
$ws.Range("I341").Value = "b"
$ws.Range("J341").Value = "Acknowledge (Backchannel)"
$ws.Range("I358").Value = "aa"
$ws.Range("J358").Value = "Agree/Accept"
$ws.Range("I363").Value = "b"
$ws.Range("J363").Value = "Acknowledge (Backchannel)"
$ws.Range("I384").Value = "sv"
$ws.Range("J384").Value = "Statement-opinion"
$ws.Range("I388").Value = "sd"
$ws.Range("J388").Value = "Statement-non-opinion"
$ws.Range("I402").Value = "b"
$ws.Range("J402").Value = "Acknowledge (Backchannel)"
$ws.Range("I403").Value = "b"
$ws.Range("J403").Value = "Acknowledge (Backchannel)"
$ws.Range("I420").Value = "sv"
$ws.Range("J420").Value = "Statement-opinion"
$ws.Range("I429").Value = "%"
$ws.Range("J429").Value = "Uninterpretable"
$ws.Range("I430").Value = "sv"
$ws.Range("J430").Value = "Statement-opinion"
$ws.Range("I433").Value = "b"
$ws.Range("J433").Value = "Acknowledge (Backchannel)"
$ws.Range("I434").Value = "b"
$ws.Range("J434").Value = "Acknowledge (Backchannel)"
$ws.Range("I438").Value = "b"
$ws.Range("J438").Value = "Acknowledge (Backchannel)"
